$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number (e.g. "1.005", "80.82", "106.40")
# must be forced to Text format first so Excel keeps them as strings, matching the
# original inlineStr cells (otherwise COM auto-converts them to numeric values and
# things like trailing zeros / exact formatting would be lost).

$ws.Range("D2").Value = '27.978.04'
$ws.Range("E2").Value = '  -3.73%  '
$ws.Range("D3").Value = '1.868.06'
$ws.Range("E3").Value = '  -2.85%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  -0.25%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '318.61'
$ws.Range("E5").Value = '  -2.30%  '
$ws.Range("E6").Value = '  -0.26%  '
$ws.Range("E7").Value = '  -5.79%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3737'
$ws.Range("E8").Value = '  -2.19%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07451'
$ws.Range("E9").Value = '  -4.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9317'
$ws.Range("E10").Value = '  -4.78%  '
$ws.Range("E11").Value = '  -6.16%  '
$ws.Range("D12").Value = '1.944.95'
$ws.Range("E12").Value = '  +1.59%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.727'
$ws.Range("E13").Value = '  -3.37%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.429'
$ws.Range("E14").Value = '  -4.57%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06856'
$ws.Range("E15").Value = '  -2.98%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.006'
$ws.Range("E16").Value = '  -0.35%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '80.82'
$ws.Range("E17").Value = '  -4.14%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009036'
$ws.Range("E18").Value = '  -5.06%  '
$ws.Range("E19").Value = '  -0.28%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.75'
$ws.Range("E20").Value = '  -5.91%  '
$ws.Range("D21").Value = '27.960.83'
$ws.Range("E21").Value = '  -3.87%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.121'
$ws.Range("E22").Value = '  -4.15%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.01'
$ws.Range("E23").Value = '  +0.51%  '
$ws.Range("D24").Value = '2.152.07'
$ws.Range("E24").Value = '  +0.43%  '
$ws.Range("E25").Value = '  -1.54%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '153.47'
$ws.Range("E26").Value = '  -2.85%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.51'
$ws.Range("E27").Value = '  -3.13%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.500'
$ws.Range("E28").Value = '  -2.88%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '113.36'
$ws.Range("E29").Value = '  -4.00%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.694'
$ws.Range("E30").Value = '  -7.79%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08994'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.8053'
$ws.Range("E32").Value = '  -5.94%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.786'
$ws.Range("E33").Value = '  -6.52%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.175'
$ws.Range("E34").Value = '  -5.47%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.947'
$ws.Range("E35").Value = '  -2.60%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.004'
$ws.Range("E36").Value = '  -0.18%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.05492'
$ws.Range("E37").Value = '  -3.36%  '
$ws.Range("E38").Value = '  -3.72%  '
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.043'
$ws.Range("E39").Value = '  -3.63%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01974'
$ws.Range("E40").Value = '  -3.45%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5234'
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.971'
$ws.Range("E42").Value = '  -6.98%  '
$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1684'
$ws.Range("E43").Value = '  -4.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.761'
$ws.Range("E44").Value = '  -5.87%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.06732'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4875'
$ws.Range("E46").Value = '  -6.06%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.45'
$ws.Range("E47").Value = '  -7.21%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '106.40'
$ws.Range("E48").Value = '  -3.53%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.003'
$ws.Range("E49").Value = '  -0.32%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.671'
$ws.Range("E50").Value = '  -5.35%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.879'
$ws.Range("E51").Value = '  -14.61%  '
